$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: append a new data row (columns A..G) to a sheet, reusing the
# number-formatting already established by the sheet's existing rows so the
# new rows look exactly like the historic ones:
#   - Column A holds a literal "yyyy-mm-dd" text label (must NOT be
#     auto-converted to a date serial, so the cell is pre-formatted as Text).
#   - Column B holds the numeric run-timestamp serial, formatted like the
#     existing timestamp column (copied from the prior row).
#   - Columns C..G are plain text/numbers and can just take the column's
#     default (General) formatting.
# ---------------------------------------------------------------------------
function Add-HistoryRow {
    param(
        $ws,
        [int]$row,
        [int]$formatSourceRow,
        [string]$runDate,
        [double]$runTime,
        [string]$sprintName,
        [double]$total,
        [double]$pass,
        [double]$fail,
        [double]$timeTaken
    )

    # Column A: force Text format first so the "yyyy-mm-dd" string is kept
    # verbatim instead of being reinterpreted as a date value, then restore
    # the column's normal (General) look by pasting the prior row's format
    # over it - the cached value stays text, only the display format moves.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $runDate
    $ws.Cells.Item($formatSourceRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    # Column B: copy the timestamp column's number format from the prior row
    # so the new cell renders/serialises exactly like the others.
    $ws.Cells.Item($formatSourceRow, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $runTime

    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $timeTaken
}

# ---------------------------------------------------------------------------
# AMSIN sheet: append rows 55-60 (sprints 165 + 166 interview-history data)
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-HistoryRow $wsAmsin 55 54 "2022-08-02" 44775.65902763889 "165_fstcycle"  89 89 0 2.23
Add-HistoryRow $wsAmsin 56 55 "2022-08-02" 44775.66594167824 "165_fstcycle"  89 89 0 2.23
Add-HistoryRow $wsAmsin 57 56 "2022-08-03" 44776.69232092593 "165_scndcycle" 89 89 0 2.39
Add-HistoryRow $wsAmsin 58 57 "2022-08-04" 44777.39471997685 "165_finalrun"  89 89 0 2.31
Add-HistoryRow $wsAmsin 59 58 "2022-08-22" 44795.6746497338  "166fstcycle"   89 89 0 2.26
Add-HistoryRow $wsAmsin 60 59 "2022-08-23" 44796.90907802084 "166cyclescnd"  89 89 0 2.2

# ---------------------------------------------------------------------------
# BETA sheet: append rows 24-25 (sprints 165 + 166 beta runs)
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Add-HistoryRow $wsBeta 24 23 "2022-08-04" 44777.56463101852 "165beta"  89 89 0 2.17
Add-HistoryRow $wsBeta 25 24 "2022-08-24" 44797.54364725694 "166_beta" 89 89 0 2.25

# ---------------------------------------------------------------------------
# AMS sheet: row 29's run-timestamp gets a tiny precision correction, then
# rows 30-31 are appended (sprints 165 + 166 live runs).
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(29, 2).Value = 44756.82411372685

Add-HistoryRow $wsAms 30 28 "2022-08-04" 44777.82010793981 "165_live" 89 89 0 2.22
Add-HistoryRow $wsAms 31 30 "2022-08-24" 44797.92527305376 "166_live" 89 89 0 2.24
